$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 154; this shifts the existing rows 154:226
# down to 155:227 (values, formats and styles move with them), matching
# every per-row value change described in the diff (each old row N's data
# now lives at row N+1, with the final old row 226 ending up at row 227).
$ws.Rows("154:154").Insert()

# Populate the freshly inserted row 154 with the new weekly price record.
$ws.Range("A154").Value = 8
$ws.Range("B154").Value = "Terminal La Palmera de La Serena"
$ws.Range("C154").Value = "Coquimbo"
$ws.Range("D154").Value = 44452
$ws.Range("D154").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E154").Value = 4
$ws.Range("F154").Value = 100114001
$ws.Range("G154").Value = "Papa"
$ws.Range("H154").Value = "Cardinal"
$ws.Range("I154").Value = "1a (cosecha)"
$ws.Range("J154").Value = 3000
$ws.Range("K154").Value = 11500
$ws.Range("L154").Value = 12000
$ws.Range("M154").Value = 11750
$ws.Range("N154").Value = "`$/saco 25 kilos"
$ws.Range("O154").Value = "Provincia del Elquí"
$ws.Range("P154").Value = 470
$ws.Range("Q154").Value = 25
$ws.Range("R154").Value = "Hortaliza"
